# Auto-generated edit script applying the diff to Maduin_Profits workbook
# Updates currentAveragePrice / Leve price / profit columns (H,I,J,K,L,M,N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 4000
$ws.Range("J9").Value = 4000
$ws.Range("L9").Value = 4000
$ws.Range("N9").Value = -4338
$ws.Range("H70").Value = 2880
$ws.Range("I70").Value = 2666.6667
$ws.Range("K70").Value = 8000.000100000001
$ws.Range("M70").Value = -7730.000100000001
$ws.Range("H73").Value = 2880
$ws.Range("I73").Value = 2666.6667
$ws.Range("K73").Value = 8000.000100000001
$ws.Range("M73").Value = -7064.000100000001
$ws.Range("H93").Value = 44999
$ws.Range("J93").Value = 44999
$ws.Range("L93").Value = 44999
$ws.Range("N93").Value = -49991

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 4666.6665
$ws.Range("I22").Value = 4000
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 4000
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -3701
$ws.Range("N22").Value = -5598
$ws.Range("H45").Value = 4258.7144
$ws.Range("I45").Value = 1932.3334
$ws.Range("K45").Value = 1932.3334
$ws.Range("M45").Value = -1555.3334
$ws.Range("H50").Value = 43899
$ws.Range("J50").Value = 43899
$ws.Range("L50").Value = 43899
$ws.Range("N50").Value = -45327
$ws.Range("H61").Value = 3830.125
$ws.Range("I61").Value = 3591.7144
$ws.Range("K61").Value = 3591.7144
$ws.Range("M61").Value = -3379.7144
$ws.Range("H74").Value = 1000
$ws.Range("I74").Value = 1000
$ws.Range("K74").Value = 1000
$ws.Range("M74").Value = -126
$ws.Range("H77").Value = 1000
$ws.Range("I77").Value = 1000
$ws.Range("K77").Value = 5000
$ws.Range("M77").Value = -632
$ws.Range("H110").Value = 3277.5715
$ws.Range("I110").Value = 3538.7
$ws.Range("J110").Value = 2624.75
$ws.Range("K110").Value = 3538.7
$ws.Range("L110").Value = 2624.75
$ws.Range("M110").Value = -1493.7
$ws.Range("N110").Value = -6714.75
$ws.Range("H132").Value = 1199.5
$ws.Range("I132").Value = 1199.5
$ws.Range("K132").Value = 3598.5
$ws.Range("M132").Value = -1068.5
$ws.Range("H136").Value = 3830.125
$ws.Range("I136").Value = 3591.7144
$ws.Range("K136").Value = 10775.1432
$ws.Range("M136").Value = -8225.143199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1704
$ws.Range("I86").Value = 2054.8572
$ws.Range("J86").Value = 885.3333
$ws.Range("K86").Value = 2054.8572
$ws.Range("L86").Value = 885.3333
$ws.Range("M86").Value = -931.8571999999999
$ws.Range("N86").Value = -3131.3333
$ws.Range("H89").Value = 1704
$ws.Range("I89").Value = 2054.8572
$ws.Range("J89").Value = 885.3333
$ws.Range("K89").Value = 10274.286
$ws.Range("L89").Value = 4426.6665
$ws.Range("M89").Value = -4658.286
$ws.Range("N89").Value = -15658.6665
$ws.Range("H105").Value = 4682
$ws.Range("I105").Value = 4234
$ws.Range("K105").Value = 4234
$ws.Range("M105").Value = -2487
$ws.Range("H107").Value = 676.3333
$ws.Range("I107").Value = 554.6667
$ws.Range("J107").Value = 798
$ws.Range("K107").Value = 554.6667
$ws.Range("L107").Value = 798
$ws.Range("M107").Value = 1365.3333
$ws.Range("N107").Value = -4638
$ws.Range("H134").Value = 8000
$ws.Range("I134").Value = 9333.333000000001
$ws.Range("K134").Value = 27999.999
$ws.Range("M134").Value = -25464.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 17428742
$ws.Range("I86").Value = 17428742
$ws.Range("K86").Value = 17428742
$ws.Range("M86").Value = -17427619
$ws.Range("H89").Value = 17428742
$ws.Range("I89").Value = 17428742
$ws.Range("K89").Value = 87143710
$ws.Range("M89").Value = -87138094

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 135
$ws.Range("I7").Value = 140.5
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 421.5
$ws.Range("L7").Value = 240
$ws.Range("M7").Value = -309.5
$ws.Range("N7").Value = -464
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11626
$ws.Range("I70").Value = 11626
$ws.Range("K70").Value = 11626
$ws.Range("M70").Value = -11356
$ws.Range("H73").Value = 11626
$ws.Range("I73").Value = 11626
$ws.Range("K73").Value = 11626
$ws.Range("M73").Value = -10690
$ws.Range("H107").Value = 739.875
$ws.Range("I107").Value = 511.5
$ws.Range("J107").Value = 1425
$ws.Range("K107").Value = 511.5
$ws.Range("L107").Value = 1425
$ws.Range("M107").Value = 1408.5
$ws.Range("N107").Value = -5265
$ws.Range("H113").Value = 776.6667
$ws.Range("J113").Value = 840
$ws.Range("L113").Value = 840
$ws.Range("N113").Value = -5180
$ws.Range("H132").Value = 4597.8335
$ws.Range("I132").Value = 4402.25
$ws.Range("K132").Value = 13206.75
$ws.Range("M132").Value = -10676.75
$ws.Range("I133").Value = 49992
$ws.Range("J133").Value = 49996
$ws.Range("K133").Value = 49992
$ws.Range("L133").Value = 49996
$ws.Range("M133").Value = -44932
$ws.Range("N133").Value = -60116

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2081.6
$ws.Range("I7").Value = 1832.6923
$ws.Range("K7").Value = 1832.6923
$ws.Range("M7").Value = -1720.6923
$ws.Range("H122").Value = 5062.4
$ws.Range("I122").Value = 5162.1665
$ws.Range("J122").Value = 4663.3335
$ws.Range("K122").Value = 15486.4995
$ws.Range("L122").Value = 13990.0005
$ws.Range("M122").Value = -13036.4995
$ws.Range("N122").Value = -18890.0005
$ws.Range("H126").Value = 2081.6
$ws.Range("I126").Value = 1832.6923
$ws.Range("K126").Value = 5498.0769
$ws.Range("M126").Value = -3028.0769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 15500
$ws.Range("J39").Value = 15500
$ws.Range("L39").Value = 15500
$ws.Range("N39").Value = -16326
$ws.Range("H107").Value = 447.2857
$ws.Range("I107").Value = 464.66666
$ws.Range("K107").Value = 1393.99998
$ws.Range("M107").Value = 526.0000199999999
$ws.Range("H126").Value = 975
$ws.Range("I126").Value = 1065.75
$ws.Range("J126").Value = 733
$ws.Range("K126").Value = 3197.25
$ws.Range("L126").Value = 2199
$ws.Range("M126").Value = -727.25
$ws.Range("N126").Value = -7139
